# Auto-generated Excel COM-interop script to apply scheduled runner updates
# to the Zalera_Profits workbook (per-sheet leve profit recalculations).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1662.6666
$ws.Range("I9").Value = 214.72223
$ws.Range("K9").Value = 214.72223
$ws.Range("M9").Value = -45.72223
$ws.Range("H15").Value = 2948102.8
$ws.Range("I15").Value = 2948102.8
$ws.Range("K15").Value = 8844308.399999999
$ws.Range("M15").Value = -8844139.399999999
$ws.Range("H101").Value = 1158.8462
$ws.Range("J101").Value = 494
$ws.Range("L101").Value = 1482
$ws.Range("N101").Value = -4726
$ws.Range("H112").Value = 2922.3958
$ws.Range("J112").Value = 2922.3958
$ws.Range("L112").Value = 8767.187399999999
$ws.Range("N112").Value = -10983.1874
$ws.Range("H137").Value = 11371398
$ws.Range("J137").Value = 7677.1333
$ws.Range("L137").Value = 23031.3999
$ws.Range("N137").Value = -28131.3999
$ws.Range("H138").Value = 2477.625
$ws.Range("I138").Value = 1276.7407
$ws.Range("J138").Value = 3353.946
$ws.Range("K138").Value = 3830.2221
$ws.Range("L138").Value = 10061.838
$ws.Range("M138").Value = 1309.7779
$ws.Range("N138").Value = -20341.838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22372.375
$ws.Range("I32").Value = 22210.086
$ws.Range("K32").Value = 22210.086
$ws.Range("M32").Value = -21923.086
$ws.Range("H45").Value = 6907.9473
$ws.Range("I45").Value = 8987.462
$ws.Range("J45").Value = 2402.3333
$ws.Range("K45").Value = 8987.462
$ws.Range("L45").Value = 2402.3333
$ws.Range("M45").Value = -8610.462
$ws.Range("N45").Value = -3156.3333
$ws.Range("H46").Value = 11883.333
$ws.Range("I46").Value = 7999
$ws.Range("J46").Value = 13825.5
$ws.Range("K46").Value = 7999
$ws.Range("L46").Value = 13825.5
$ws.Range("M46").Value = -7680
$ws.Range("N46").Value = -14463.5
$ws.Range("H61").Value = 9731.727999999999
$ws.Range("I61").Value = 7800
$ws.Range("J61").Value = 10456.125
$ws.Range("K61").Value = 7800
$ws.Range("L61").Value = 10456.125
$ws.Range("M61").Value = -7588
$ws.Range("N61").Value = -10880.125
$ws.Range("H122").Value = 5093.2
$ws.Range("I122").Value = 5116.75
$ws.Range("K122").Value = 15350.25
$ws.Range("M122").Value = -12900.25
$ws.Range("H132").Value = 3497.3914
$ws.Range("I132").Value = 2312.7368
$ws.Range("K132").Value = 6938.2104
$ws.Range("M132").Value = -4408.2104
$ws.Range("H136").Value = 9731.727999999999
$ws.Range("I136").Value = 7800
$ws.Range("J136").Value = 10456.125
$ws.Range("K136").Value = 23400
$ws.Range("L136").Value = 31368.375
$ws.Range("M136").Value = -20850
$ws.Range("N136").Value = -36468.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 8000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null
$ws.Range("H105").Value = 41678704
$ws.Range("I105").Value = 50013924
$ws.Range("K105").Value = 50013924
$ws.Range("M105").Value = -50012177
$ws.Range("H134").Value = 6555.85
$ws.Range("I134").Value = 2816
$ws.Range("K134").Value = 8448
$ws.Range("M134").Value = -5913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33338680
$ws.Range("I31").Value = 100001190
$ws.Range("J31").Value = 7425.45
$ws.Range("K31").Value = 100001190
$ws.Range("L31").Value = 7425.45
$ws.Range("M31").Value = -100000895
$ws.Range("N31").Value = -8015.45
$ws.Range("H34").Value = 33338680
$ws.Range("I34").Value = 100001190
$ws.Range("J34").Value = 7425.45
$ws.Range("K34").Value = 100001190
$ws.Range("L34").Value = 7425.45
$ws.Range("M34").Value = -100000988
$ws.Range("N34").Value = -7829.45
$ws.Range("H58").Value = 5611.593
$ws.Range("I58").Value = 3662.7058
$ws.Range("K58").Value = 3662.7058
$ws.Range("M58").Value = -3459.7058
$ws.Range("H122").Value = 72693.28999999999
$ws.Range("I122").Value = 91883.37
$ws.Range("K122").Value = 275650.11
$ws.Range("M122").Value = -273200.11
$ws.Range("H134").Value = 5110.3335
$ws.Range("I134").Value = 4663.207
$ws.Range("K134").Value = 13989.621
$ws.Range("M134").Value = -11454.621
$ws.Range("H136").Value = 5611.593
$ws.Range("I136").Value = 3662.7058
$ws.Range("K136").Value = 10988.1174
$ws.Range("M136").Value = -8438.117400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1413938
$ws.Range("I4").Value = 2296359
$ws.Range("J4").Value = 64353.176
$ws.Range("K4").Value = 6889077
$ws.Range("L4").Value = 193059.528
$ws.Range("M4").Value = -6888965
$ws.Range("N4").Value = -193283.528
$ws.Range("H5").Value = 859.0769
$ws.Range("I5").Value = 579.6
$ws.Range("J5").Value = 1033.75
$ws.Range("K5").Value = 1738.8
$ws.Range("L5").Value = 3101.25
$ws.Range("M5").Value = -1626.8
$ws.Range("N5").Value = -3325.25
$ws.Range("H29").Value = 677.1
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 821.375
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 2464.125
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -3018.125
$ws.Range("H32").Value = 947.9
$ws.Range("I32").Value = 630
$ws.Range("J32").Value = 1159.8334
$ws.Range("K32").Value = 1890
$ws.Range("L32").Value = 3479.5002
$ws.Range("M32").Value = -1607
$ws.Range("N32").Value = -4045.5002
$ws.Range("H46").Value = 153
$ws.Range("J46").Value = 153
$ws.Range("L46").Value = 459
$ws.Range("N46").Value = -641
$ws.Range("H133").Value = 2329.8
$ws.Range("I133").Value = 2329.8
$ws.Range("K133").Value = 6989.400000000001
$ws.Range("M133").Value = -1929.400000000001
$ws.Range("H135").Value = 859.0769
$ws.Range("I135").Value = 579.6
$ws.Range("J135").Value = 1033.75
$ws.Range("K135").Value = 5216.400000000001
$ws.Range("L135").Value = 9303.75
$ws.Range("M135").Value = -2681.400000000001
$ws.Range("N135").Value = -14373.75
$ws.Range("H136").Value = 2428.3
$ws.Range("I136").Value = 454.85715
$ws.Range("K136").Value = 1364.57145
$ws.Range("M136").Value = 3735.42855
$ws.Range("H139").Value = 50001430
$ws.Range("I139").Value = 50001430
$ws.Range("K139").Value = 150004290
$ws.Range("M139").Value = -149999150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 750.3158
$ws.Range("I2").Value = 1106.4166
$ws.Range("K2").Value = 1106.4166
$ws.Range("M2").Value = -993.4166
$ws.Range("H122").Value = 7953.778
$ws.Range("I122").Value = 8009.9414
$ws.Range("K122").Value = 24029.8242
$ws.Range("M122").Value = -21579.8242
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 26500
$ws.Range("J43").Value = 26500
$ws.Range("L43").Value = 26500
$ws.Range("N43").Value = -26886
$ws.Range("I46").Value = 1600.5
$ws.Range("J46").Value = 4726.3955
$ws.Range("K46").Value = 1600.5
$ws.Range("L46").Value = 4726.3955
$ws.Range("M46").Value = -1412.5
$ws.Range("N46").Value = -5102.3955
$ws.Range("H61").Value = 1514.2727
$ws.Range("I61").Value = 1141.5
$ws.Range("K61").Value = 1141.5
$ws.Range("M61").Value = -939.5
$ws.Range("H68").Value = 10911.444
$ws.Range("J68").Value = 36449.5
$ws.Range("L68").Value = 36449.5
$ws.Range("N68").Value = -37947.5
$ws.Range("H71").Value = 10911.444
$ws.Range("J71").Value = 36449.5
$ws.Range("L71").Value = 182247.5
$ws.Range("N71").Value = -189735.5
$ws.Range("H82").Value = 3295.889
$ws.Range("J82").Value = 4332.6
$ws.Range("L82").Value = 4332.6
$ws.Range("N82").Value = -5054.6
$ws.Range("H85").Value = 3295.889
$ws.Range("J85").Value = 4332.6
$ws.Range("L85").Value = 4332.6
$ws.Range("N85").Value = -6828.6
$ws.Range("H94").Value = 53000
$ws.Range("J94").Value = 53000
$ws.Range("L94").Value = 53000
$ws.Range("N94").Value = -54352
$ws.Range("H113").Value = 1514.2727
$ws.Range("I113").Value = 1141.5
$ws.Range("K113").Value = 1141.5
$ws.Range("M113").Value = 1028.5
$ws.Range("H122").Value = 7288.4546
$ws.Range("I122").Value = 7317.3
$ws.Range("K122").Value = 21951.9
$ws.Range("M122").Value = -19501.9
$ws.Range("H136").Value = 5111.6284
$ws.Range("I136").Value = 3055.6667
$ws.Range("K136").Value = 9167.000100000001
$ws.Range("M136").Value = -6617.000100000001
$ws.Range("H140").Value = 63550.57
$ws.Range("J140").Value = 63550.57
$ws.Range("L140").Value = 63550.57
$ws.Range("N140").Value = -73910.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26991.666
$ws.Range("I62").Value = 26991.666
$ws.Range("K62").Value = 26991.666
$ws.Range("M62").Value = -26367.666
$ws.Range("H65").Value = 26991.666
$ws.Range("I65").Value = 26991.666
$ws.Range("K65").Value = 134958.33
$ws.Range("M65").Value = -131838.33
$ws.Range("H100").Value = 2678.3572
$ws.Range("I100").Value = 3076.6365
$ws.Range("K100").Value = 6153.273
$ws.Range("M100").Value = -5612.273
$ws.Range("H122").Value = 4999.5
$ws.Range("I122").Value = 4999.5
$ws.Range("K122").Value = 14998.5
$ws.Range("M122").Value = -12548.5
$ws.Range("H132").Value = 4226.1177
$ws.Range("I132").Value = 2442.7222
$ws.Range("K132").Value = 7328.1666
$ws.Range("M132").Value = -4798.1666

